$d = $word.ActiveDocument

# The paragraph that credits the astronomical star-map authors and links
# to the yearly map archive. We replace its whole run-run-run content
# (built out of many differently-formatted / spell-check-wrapped runs)
# with a single plain run holding the full sentence, bumping the year
# in the URL from 2018 to 2022.
$oldSentence = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newSentence = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Os mapas de*GaNight*") {
        $start = $p.Range.Start
        $end = $p.Range.End
        # Exclude the trailing paragraph mark from the replaced range.
        $r = $d.Range($start, $end - 1)
        $r.Text = ""
        $r.InsertAfter($newSentence)
        break
    }
}
